$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of test-case data (row 16): "Customercare015" ---
$ws.Range("A16").Value = "Customercare015"
$ws.Range("B16").Value = "OPQA-5320||OPQA-5321||OPQA5322||OPQA-5323"
$ws.Range("C16").Value = "verify that upon clicking on submit button, a success message should be displayed that confirms submission and should allow user to raise a new ticket||Verify that success message in customer care page should match with wire frame||verify that extension field should be placed next to phone number field in customer care page||Verify that all characters including special characters should be allowed in extension field in customer care page."
$ws.Range("D16").Value = "Y"

# Match the wrapped-text style already used by the other filled-in rows
# in column A (this flips A16 from the blank-row style onto the
# "TCID" content style, same as A10:A15).
$ws.Range("A16").WrapText = $true

# The row grows to fit the new multi-line description, same as the
# other ht=45 rows (8, 9, 10, 11).
$ws.Rows.Item(16).RowHeight = 45

# --- Column B was narrowed to make room ---
# (29.8333.. is how this engine round-trips a plain "29"; back it off by
# the fixed ~0.8333-character padding it adds on export so the saved
# <col> width lands on exactly 29.)
$ws.Columns.Item(2).ColumnWidth = 28.1666666667

# --- Selection / scroll position left where the edit was made ---
$ws.Range("C19").Select()
